$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell for column I, formatted like the existing H1 header
$ws.Range("I1").Value = "Total Clan Stars"
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)

# Values for rows 2-16 (Total Clan Stars), each styled like the row's H cell
$values = @{
    2  = 37
    3  = 37
    4  = 37
    5  = 37
    6  = 37
    7  = 37
    8  = 37
    9  = 37
    10 = 37
    11 = 37
    12 = 0
    13 = 0
    14 = 0
    15 = 0
    16 = 0
}

foreach ($row in 2..16) {
    $ws.Cells.Item($row, 9).Value = $values[$row]
    $ws.Cells.Item($row, 8).Copy()
    $ws.Cells.Item($row, 9).PasteSpecial(-4122)
}

$excel.CutCopyMode = 0
